$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new quarterly data point was added at the top of the error table (row 2),
# so every existing row shifts down by one (row 2 -> row 3, ... row 10 -> row 11),
# and the previous row 11 (oldest window entry) drops off the bottom.
for ($r = 11; $r -ge 3; $r--) {
    $src = $r - 1
    $v = $ws.Range("B$src`:G$src").Value2
    $ws.Range("B$r`:G$r").Value = $v
}

# Write the newly computed error statistics for the new top row (Q6).
$ws.Range("B2").Value = 0.01792803175822741
$ws.Range("C2").Value = 0.6873742451724051
$ws.Range("D2").Value = 1.225852550159523
$ws.Range("E2").Value = 1.107182256974669
$ws.Range("F2").Value = 1.137372497882221
$ws.Range("G2").Value = 19
